# "Add files via upload" — log entries for Dec 4-5 2023 (rows 121-137 on Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Total Time for the Dec 4 2023 block (header row 121)
$ws.Range("D121").Value = 9

# Dec 4 2023 — afternoon slots (rows 123-132)
$ws.Range("B123").Value = "LED Ring Layout design"
$ws.Range("C123").Value = "design"

$ws.Range("B124").Value = "LED Ring Layout design"
$ws.Range("C124").Value = "design"

$ws.Range("B125").Value = "LED Ring Layout design"
$ws.Range("C125").Value = "design"

$ws.Range("B126").Value = "LED Ring Layout design"
$ws.Range("C126").Value = "design"

$ws.Range("B127").Value = "search for customers"
$ws.Range("C127").Value = "Support"

$ws.Range("B128").Value = "search for customers"
$ws.Range("C128").Value = "Support"

$ws.Range("B129").Value = "search for customers"
$ws.Range("C129").Value = "Support"

$ws.Range("B130").Value = "LED Ring Layout design"
$ws.Range("C130").Value = "design"

$ws.Range("B131").Value = "LED Ring Layout design"
$ws.Range("C131").Value = "design"

$ws.Range("B132").Value = "LED Ring Layout design"
$ws.Range("C132").Value = "design"

# Dec 5 2023 — morning slots (rows 135-137)
$ws.Range("B135").Value = "LED Ring Layout design"
$ws.Range("C135").Value = "design"

$ws.Range("B136").Value = "LED Ring Layout design"
$ws.Range("C136").Value = "design"

$ws.Range("B137").Value = "LED Ring Layout UPLOAD"
$ws.Range("C137").Value = "design"

# restore the selection/active cell left by the author when they saved
$ws.Range("D138").Select()
